$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds a "Date" header in row 1 and a per-row game-date string in
# every row below it. The dates were originally written in an ambiguous
# "M-D-YYYY-YY" form, e.g. "6-23-2012-13", which (per the commit message)
# was off by a day because of the way the NBA stats were shown. Replace
# every occurrence of that bad string with the corrected ISO (YYYY-MM-DD)
# date, "2013-06-23", keeping the cell as plain text so Excel does not
# reinterpret the corrected value as a serial date.
$oldValue = "6-23-2012-13"
$newValue = "2013-06-23"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    if ($cell.Value2 -eq $oldValue) {
        $cell.NumberFormat = "@"
        $cell.Value2 = $newValue
    }
}
